$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 267
$ws.Range("B267").Value = 6707529
$ws.Range("F267").Value = "Internacional"
$ws.Range("G267").Value = "Bragantino"
$ws.Range("H267").Value = 1
$ws.Range("I267").Value = 0
$ws.Range("J267").Value = "H"
$ws.Range("K267").Value = 2.6
$ws.Range("L267").Value = 3.3
$ws.Range("M267").Value = 2.625
$ws.Range("N267").Value = 2.4
$ws.Range("O267").Value = 3.4
$ws.Range("P267").Value = 2.9
$ws.Range("Q267").Value = -0.25
$ws.Range("R267").Value = 2.05
$ws.Range("S267").Value = 1.8
$ws.Range("T267").Value = 2.25
$ws.Range("U267").Value = 1.875
$ws.Range("V267").Value = 1.975
$ws.Range("W267").Value = 1.4
$ws.Range("X267").Value = -1
$ws.Range("Y267").Value = -1
$ws.Range("Z267").Value = 1.05
$ws.Range("AA267").Value = -1
$ws.Range("AB267").Value = -1
$ws.Range("AC267").Value = 0.9750000000000001

# Row 268
$ws.Range("B268").Value = 6637300
$ws.Range("F268").Value = "Sao Paulo"
$ws.Range("G268").Value = "Cuiaba"
$ws.Range("H268").Value = 0
$ws.Range("I268").Value = 0
$ws.Range("J268").Value = "D"
$ws.Range("K268").Value = 1.727
$ws.Range("L268").Value = 3.4
$ws.Range("M268").Value = 5.25
$ws.Range("N268").Value = 1.615
$ws.Range("O268").Value = 3.6
$ws.Range("P268").Value = 7
$ws.Range("Q268").Value = -0.75
$ws.Range("R268").Value = 1.825
$ws.Range("S268").Value = 2.025
$ws.Range("T268").Value = 2
$ws.Range("U268").Value = 1.925
$ws.Range("V268").Value = 1.925
$ws.Range("W268").Value = -1
$ws.Range("X268").Value = 2.6
$ws.Range("Y268").Value = -1
$ws.Range("Z268").Value = -1
$ws.Range("AA268").Value = 1.025
$ws.Range("AB268").Value = -1
$ws.Range("AC268").Value = 0.925

# Row 269
$ws.Range("B269").Value = 6704082
$ws.Range("F269").Value = "Fortaleza EC"
$ws.Range("G269").Value = "Palmeiras"
$ws.Range("H269").Value = 2
$ws.Range("I269").Value = 2
$ws.Range("J269").Value = "D"
$ws.Range("K269").Value = 3.25
$ws.Range("L269").Value = 3.2
$ws.Range("M269").Value = 2.25
$ws.Range("N269").Value = 3.1
$ws.Range("O269").Value = 3.1
$ws.Range("P269").Value = 2.375
$ws.Range("Q269").Value = 0.25
$ws.Range("R269").Value = 1.8
$ws.Range("S269").Value = 2.05
$ws.Range("T269").Value = 2.25
$ws.Range("U269").Value = 2.025
$ws.Range("V269").Value = 1.825
$ws.Range("W269").Value = -1
$ws.Range("X269").Value = 2.1
$ws.Range("Y269").Value = -1
$ws.Range("Z269").Value = 0.4
$ws.Range("AA269").Value = -0.5
$ws.Range("AB269").Value = 1.025
$ws.Range("AC269").Value = -1

# Row 270
$ws.Range("B270").Value = 6704396
$ws.Range("F270").Value = "America MG"
$ws.Range("G270").Value = "Flamengo"
$ws.Range("H270").Value = 0
$ws.Range("I270").Value = 3
$ws.Range("J270").Value = "A"
$ws.Range("K270").Value = 6.5
$ws.Range("L270").Value = 3.75
$ws.Range("M270").Value = 1.533
$ws.Range("N270").Value = 7.5
$ws.Range("O270").Value = 4.75
$ws.Range("P270").Value = 1.4
$ws.Range("Q270").Value = 1.25
$ws.Range("R270").Value = 1.95
$ws.Range("S270").Value = 1.9
$ws.Range("T270").Value = 2.75
$ws.Range("U270").Value = 1.8
$ws.Range("V270").Value = 2.05
$ws.Range("W270").Value = -1
$ws.Range("X270").Value = -1
$ws.Range("Y270").Value = 0.3999999999999999
$ws.Range("Z270").Value = -1
$ws.Range("AA270").Value = 0.8999999999999999
$ws.Range("AB270").Value = 0.4
$ws.Range("AC270").Value = -0.5

# Row 286
$ws.Range("B286").Value = 6677292
$ws.Range("F286").Value = "Bragantino"
$ws.Range("G286").Value = "Coritiba"
$ws.Range("H286").Value = 1
$ws.Range("I286").Value = 0
$ws.Range("J286").Value = "H"
$ws.Range("K286").Value = 1.3
$ws.Range("L286").Value = 5.5
$ws.Range("M286").Value = 9
$ws.Range("N286").Value = 1.2
$ws.Range("O286").Value = 7.5
$ws.Range("P286").Value = 12
$ws.Range("Q286").Value = -2
$ws.Range("R286").Value = 1.98
$ws.Range("S286").Value = 1.92
$ws.Range("T286").Value = 3.25
$ws.Range("U286").Value = 1.925
$ws.Range("V286").Value = 1.925
$ws.Range("W286").Value = 0.2
$ws.Range("X286").Value = -1
$ws.Range("Y286").Value = -1
$ws.Range("Z286").Value = -1
$ws.Range("AA286").Value = 0.9199999999999999
$ws.Range("AB286").Value = -1
$ws.Range("AC286").Value = 0.925

# Row 287
$ws.Range("B287").Value = 6677290
$ws.Range("F287").Value = "Gremio"
$ws.Range("G287").Value = "Vasco da Gama"
$ws.Range("H287").Value = 1
$ws.Range("I287").Value = 0
$ws.Range("J287").Value = "H"
$ws.Range("K287").Value = 1.727
$ws.Range("L287").Value = 3.6
$ws.Range("M287").Value = 4.75
$ws.Range("N287").Value = 1.85
$ws.Range("O287").Value = 3.5
$ws.Range("P287").Value = 4.5
$ws.Range("Q287").Value = -0.5
$ws.Range("R287").Value = 1.825
$ws.Range("S287").Value = 2.025
$ws.Range("T287").Value = 2.5
$ws.Range("U287").Value = 1.9
$ws.Range("V287").Value = 1.95
$ws.Range("W287").Value = 0.8500000000000001
$ws.Range("X287").Value = -1
$ws.Range("Y287").Value = -1
$ws.Range("Z287").Value = 0.825
$ws.Range("AA287").Value = -1
$ws.Range("AB287").Value = -1
$ws.Range("AC287").Value = 0.95

# Row 288
$ws.Range("B288").Value = 6733768
$ws.Range("F288").Value = "Athletico Paranaense"
$ws.Range("G288").Value = "Santos"
$ws.Range("H288").Value = 3
$ws.Range("I288").Value = 0
$ws.Range("J288").Value = "H"
$ws.Range("K288").Value = 1.85
$ws.Range("L288").Value = 3.4
$ws.Range("M288").Value = 4.333
$ws.Range("N288").Value = 1.95
$ws.Range("O288").Value = 3.25
$ws.Range("P288").Value = 4.2
$ws.Range("Q288").Value = -0.5
$ws.Range("R288").Value = 1.975
$ws.Range("S288").Value = 1.875
$ws.Range("T288").Value = 2.25
$ws.Range("U288").Value = 1.875
$ws.Range("V288").Value = 1.975
$ws.Range("W288").Value = 0.95
$ws.Range("X288").Value = -1
$ws.Range("Y288").Value = -1
$ws.Range("Z288").Value = 0.9750000000000001
$ws.Range("AA288").Value = -1
$ws.Range("AB288").Value = 0.875
$ws.Range("AC288").Value = -1

# Row 289
$ws.Range("B289").Value = 6733112
$ws.Range("F289").Value = "Fortaleza EC"
$ws.Range("G289").Value = "Goias"
$ws.Range("H289").Value = 1
$ws.Range("I289").Value = 0
$ws.Range("J289").Value = "H"
$ws.Range("K289").Value = 1.75
$ws.Range("L289").Value = 3.5
$ws.Range("M289").Value = 4.75
$ws.Range("N289").Value = 1.4
$ws.Range("O289").Value = 4.333
$ws.Range("P289").Value = 9
$ws.Range("Q289").Value = -1.25
$ws.Range("R289").Value = 1.96
$ws.Range("S289").Value = 1.94
$ws.Range("T289").Value = 2.5
$ws.Range("U289").Value = 1.825
$ws.Range("V289").Value = 2.025
$ws.Range("W289").Value = 0.3999999999999999
$ws.Range("X289").Value = -1
$ws.Range("Y289").Value = -1
$ws.Range("Z289").Value = -0.5
$ws.Range("AA289").Value = 0.47
$ws.Range("AB289").Value = -1
$ws.Range("AC289").Value = 1.025

# Row 290
$ws.Range("B290").Value = 6705049
$ws.Range("F290").Value = "Botafogo"
$ws.Range("G290").Value = "Cruzeiro"
$ws.Range("H290").Value = 0
$ws.Range("I290").Value = 0
$ws.Range("J290").Value = "D"
$ws.Range("K290").Value = 1.833
$ws.Range("L290").Value = 3.4
$ws.Range("M290").Value = 4.5
$ws.Range("N290").Value = 1.909
$ws.Range("O290").Value = 3.4
$ws.Range("P290").Value = 4.2
$ws.Range("Q290").Value = -0.5
$ws.Range("R290").Value = 1.96
$ws.Range("S290").Value = 1.94
$ws.Range("T290").Value = 2.25
$ws.Range("U290").Value = 1.925
$ws.Range("V290").Value = 1.925
$ws.Range("W290").Value = -1
$ws.Range("X290").Value = 2.4
$ws.Range("Y290").Value = -1
$ws.Range("Z290").Value = -1
$ws.Range("AA290").Value = 0.9399999999999999
$ws.Range("AB290").Value = -1
$ws.Range("AC290").Value = 0.925

# Row 291
$ws.Range("B291").Value = 6707715
$ws.Range("F291").Value = "America MG"
$ws.Range("G291").Value = "EC Bahia"
$ws.Range("H291").Value = 3
$ws.Range("I291").Value = 2
$ws.Range("J291").Value = "H"
$ws.Range("K291").Value = 4
$ws.Range("L291").Value = 3.4
$ws.Range("M291").Value = 1.909
$ws.Range("N291").Value = 4
$ws.Range("O291").Value = 3.75
$ws.Range("P291").Value = 1.909
$ws.Range("Q291").Value = 0.5
$ws.Range("R291").Value = 1.9
$ws.Range("S291").Value = 1.95
$ws.Range("T291").Value = 2.75
$ws.Range("U291").Value = 2.025
$ws.Range("V291").Value = 1.825
$ws.Range("W291").Value = 3
$ws.Range("X291").Value = -1
$ws.Range("Y291").Value = -1
$ws.Range("Z291").Value = 0.8999999999999999
$ws.Range("AA291").Value = -1
$ws.Range("AB291").Value = 1.025
$ws.Range("AC291").Value = -1
